# Apply the "append the parameters for curvelet transform" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Fix up a few description strings (column D) that were reworded ---
$ws.Range("D5").Value2  = "thresh_im2 is for a hard value, main adjustable parameters"
$ws.Range("D14").Value2 = "maximum dangle angler difference at cross-link(cos(10*pi/180))"
$ws.Range("D18").Value2 = "distance for linknig same-oriented fibers"

# --- 2. Remove the stray leftover cells that aren't part of the parameter table ---
$ws.Range("O3").ClearContents()   | Out-Null
$ws.Range("J4:J5").ClearContents() | Out-Null
$ws.Range("J9").ClearContents()   | Out-Null
$ws.Range("J41").EntireRow.Delete() | Out-Null

# --- 3. Append the two new curvelet-transform parameter rows ---
$ws.Cells.Item(28, 1).Value2 = 28
$ws.Cells.Item(28, 2).Value2 = "pct"
$ws.Cells.Item(28, 3).Value2 = 0.2
$ws.Cells.Item(28, 4).Value2 = "Percentile of the remaining curvelet coeffs"

$ws.Cells.Item(29, 1).Value2 = 29
$ws.Cells.Item(29, 2).Value2 = "ss"
$ws.Cells.Item(29, 3).Value2 = 3
$ws.Cells.Item(29, 4).Value2 = "Number of selected scales"

# --- 4. Column widths for B and C narrowed slightly ---
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666

# --- 5. Update the remembered selection/cursor cell ---
$ws.Range("Q14").Select() | Out-Null

Write-Output "edit applied"
